$d = $word.ActiveDocument

# 1. Change the text of the "What are values, Variables, Operators..." (second occurrence, size 36) paragraph
$d.Content.Find.Execute(
    "What are values, Variables, Operators or function that allow me to write the code I want?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What are Operators or functions that allow me to write the code I want?", 2)

# 2. Merge the three runs " The " + "else if(){} condition statement" + " follows all the same rules as if() {}."
#    into a single run of text " The else if(){} condition statement follows all the same rules as if() {}."
#    Since Find/Replace works on text content regardless of run boundaries, replacing across the three
#    runs with the identical combined text will merge them into a single run.
$d.Content.Find.Execute(
    " The else if(){} condition statement follows all the same rules as if() {}.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " The else if(){} condition statement follows all the same rules as if() {}.", 2)

# 3. Move the _GoBack bookmark from the last paragraph (after "So you don't run into an error.")
#    to the empty paragraph right after "What are Operators or functions..." paragraph.
#    First remove the existing bookmark if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the paragraph that immediately follows the "What are Operators or functions..." paragraph.
# It is an empty paragraph with sz 44. We locate it by finding the paragraph containing our edited text,
# then taking the Next paragraph.
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t -like "*What are Operators or functions that allow me to write the code I want?*") {
        $nextPara = $para.Next()
        $d.Bookmarks.Add("_GoBack", $nextPara.Range)
        break
    }
}
